$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new "哈希" worksheet right after the existing "链表" sheet.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "哈希"

# ---------------------------------------------------------------------------
# Header row (row 1) - same headers as sheet1, bold 微软雅黑 14, centered/wrap.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "No."
$ws.Range("B1").Value = "leetcode"
$ws.Range("C1").Value = "题目"
$ws.Range("D1").Value = "解题方法"
$ws.Range("E1").Value = "解题关键词"
$ws.Range("F1").Value = "时间复杂度"
$ws.Range("G1").Value = "空间复杂度"

$hdr = $ws.Range("A1:G1")
$hdr.Font.Name = "微软雅黑"
$hdr.Font.Size = 14
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4131
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true

# Trailing empty header cells (H1:K1) that exist on sheet1 as well.
$pad = $ws.Range("H1:K1")
$pad.Font.Name = "Calibri"
$pad.Font.Size = 14
$pad.Font.Bold = $false
$pad.HorizontalAlignment = -4131
$pad.VerticalAlignment = -4108
$pad.WrapText = $true

# ---------------------------------------------------------------------------
# Data row (row 2) - the new "single number with hash / xor" entry.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 136
$ws.Range("C2").Value = "给定一个非空整数数组，除了某个元素只出现一次以外，其余每个元素均出现两次。找出那个只出现了一次的元素"
$ws.Range("E2").Value = "异或"
$ws.Range("D2").Value = "异或方案`n异或性质：`n1）一个数字和0异或【xor，^】，结果是其本身`n2）一个数字和其本身异或，结果是0`n3）异或满足交换律和结合律:a^b^b^c^a = (a^a)^(b^b)^c=0^0^c=0^c=c`n0逐个与数组中的每个数字异或操作，剩余的就是只出现一次的那个数字"
$ws.Range("F2").Value = "O(N), N是元素个数"
$ws.Range("G2").Value = "O(1)"

$dat = $ws.Range("A2:G2")
$dat.Font.Name = "微软雅黑"
$dat.Font.Size = 14
$dat.Font.Bold = $false
$dat.HorizontalAlignment = -4131
$dat.VerticalAlignment = -4108
$dat.WrapText = $true

# ---------------------------------------------------------------------------
# Row heights / column widths to match the authored layout.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 35
$ws.Rows.Item(2).RowHeight = 198

$ws.Columns.Item(1).ColumnWidth = 10.142857142857142
$ws.Columns.Item(2).ColumnWidth = 12.857142857142858
$ws.Columns.Item(3).ColumnWidth = 34.142857142857146
$ws.Columns.Item(4).ColumnWidth = 63.42857142857143
$ws.Columns.Item(5).ColumnWidth = 18.142857142857142
$ws.Columns.Item(6).ColumnWidth = 16.571428571428573
$ws.Columns.Item(7).ColumnWidth = 19.714285714285715

# ---------------------------------------------------------------------------
# Selections: sheet1 now points at F4:G4, sheet2 (new, active) points at F10.
# ---------------------------------------------------------------------------
$null = $ws1.Range("F4:G4").Select()
$null = $ws.Range("F10").Select()
